# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on the per-locale
# report sheets, as produced by a re-run of the handback report generator.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 2 corresponds to file 8f6e7f4d-... .zh-cn.xlf
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-20 18:52:08"
$wsZh.Range("H2").Value = "2016-03-20 18:52:27"

# de-de sheet: row 2 corresponds to file 8f6e7f4d-... .de-de.xlf
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-20 18:52:11"
$wsDe.Range("H2").Value = "2016-03-20 18:52:32"
